$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) and E (Volume(1h)) updates per-row. All source cells
# are stored as text (inline strings) in the original workbook, even
# when the text looks like a plain number (e.g. "596.58"), so every
# write below forces a text number-format first and restores the
# default "Normal" style afterwards to avoid leaving stray formatting.
$changes = @(
    @{ Row = 2; Col = "D"; Value = "63.103.91" },
    @{ Row = 2; Col = "E"; Value = "  -1.94%  " },
    @{ Row = 3; Col = "D"; Value = "3.130.75" },
    @{ Row = 3; Col = "E"; Value = "  -0.34%  " },
    @{ Row = 4; Col = "E"; Value = "  +0.13%  " },
    @{ Row = 5; Col = "D"; Value = "596.58" },
    @{ Row = 5; Col = "E"; Value = "  -2.33%  " },
    @{ Row = 6; Col = "D"; Value = "136.99" },
    @{ Row = 6; Col = "E"; Value = "  -4.87%  " },
    @{ Row = 7; Col = "E"; Value = "  +0.16%  " },
    @{ Row = 8; Col = "D"; Value = "3.120.89" },
    @{ Row = 8; Col = "E"; Value = "  -0.60%  " },
    @{ Row = 9; Col = "E"; Value = "  -1.87%  " },
    @{ Row = 10; Col = "D"; Value = "0.146" },
    @{ Row = 10; Col = "E"; Value = "  -3.11%  " },
    @{ Row = 11; Col = "D"; Value = "5.28" },
    @{ Row = 11; Col = "E"; Value = "  -2.16%  " },
    @{ Row = 12; Col = "D"; Value = "0.460" },
    @{ Row = 12; Col = "E"; Value = "  -3.51%  " },
    @{ Row = 13; Col = "E"; Value = "  -2.82%  " },
    @{ Row = 14; Col = "D"; Value = "34.40" },
    @{ Row = 14; Col = "E"; Value = "  -3.57%  " },
    @{ Row = 15; Col = "D"; Value = "3.643.55" },
    @{ Row = 15; Col = "E"; Value = "  -0.37%  " },
    @{ Row = 16; Col = "E"; Value = "  +1.79%  " },
    @{ Row = 17; Col = "D"; Value = "63.184.95" },
    @{ Row = 17; Col = "E"; Value = "  -1.80%  " },
    @{ Row = 18; Col = "D"; Value = "3.132.02" },
    @{ Row = 18; Col = "E"; Value = "  -0.23%  " },
    @{ Row = 19; Col = "E"; Value = "  -1.74%  " },
    @{ Row = 20; Col = "D"; Value = "477.18" },
    @{ Row = 20; Col = "E"; Value = "  -0.13%  " },
    @{ Row = 21; Col = "D"; Value = "14.22" },
    @{ Row = 21; Col = "E"; Value = "  -3.55%  " },
    @{ Row = 22; Col = "E"; Value = "  -3.23%  " },
    @{ Row = 23; Col = "D"; Value = "7.70" },
    @{ Row = 23; Col = "E"; Value = "  -1.67%  " },
    @{ Row = 24; Col = "D"; Value = "87.48" },
    @{ Row = 24; Col = "E"; Value = "  +2.63%  " },
    @{ Row = 25; Col = "D"; Value = "13.06" },
    @{ Row = 25; Col = "E"; Value = "  -4.53%  " },
    @{ Row = 26; Col = "E"; Value = "  +0.12%  " },
    @{ Row = 27; Col = "D"; Value = "2.72" },
    @{ Row = 27; Col = "E"; Value = "  -2.47%  " },
    @{ Row = 28; Col = "D"; Value = "7.21" },
    @{ Row = 28; Col = "E"; Value = "  -3.41%  " },
    @{ Row = 29; Col = "D"; Value = "8.00" },
    @{ Row = 29; Col = "E"; Value = "  -7.27%  " },
    @{ Row = 30; Col = "E"; Value = "  +0.37%  " },
    @{ Row = 31; Col = "D"; Value = "27.18" },
    @{ Row = 31; Col = "E"; Value = "  +1.75%  " },
    @{ Row = 32; Col = "E"; Value = "  +0.08%  " },
    @{ Row = 33; Col = "E"; Value = "  -7.67%  " },
    @{ Row = 34; Col = "D"; Value = "2.54" },
    @{ Row = 34; Col = "E"; Value = "  -3.70%  " },
    @{ Row = 35; Col = "E"; Value = "  -2.77%  " },
    @{ Row = 36; Col = "D"; Value = "5.84" },
    @{ Row = 36; Col = "E"; Value = "  -1.91%  " },
    @{ Row = 37; Col = "D"; Value = "51.92" },
    @{ Row = 37; Col = "E"; Value = "  -1.12%  " },
    @{ Row = 38; Col = "E"; Value = "  -4.74%  " },
    @{ Row = 39; Col = "E"; Value = "  -1.53%  " },
    @{ Row = 40; Col = "D"; Value = "423.60" },
    @{ Row = 40; Col = "E"; Value = "  -7.23%  " },
    @{ Row = 41; Col = "E"; Value = "  -0.84%  " },
    @{ Row = 42; Col = "E"; Value = "  -0.72%  " },
    @{ Row = 43; Col = "E"; Value = "  -11.16%  " },
    @{ Row = 44; Col = "D"; Value = "2.888.85" },
    @{ Row = 44; Col = "E"; Value = "  +0.81%  " },
    @{ Row = 45; Col = "D"; Value = "0.265" },
    @{ Row = 45; Col = "E"; Value = "  -0.62%  " },
    @{ Row = 48; Col = "D"; Value = "25.87" },
    @{ Row = 48; Col = "E"; Value = "  -2.67%  " },
    @{ Row = 49; Col = "E"; Value = "  -0.32%  " },
    @{ Row = 50; Col = "D"; Value = "2.29" },
    @{ Row = 50; Col = "E"; Value = "  -6.08%  " },
    @{ Row = 51; Col = "D"; Value = "118.56" },
    @{ Row = 51; Col = "E"; Value = "  -1.65%  " }
)

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.Style = "Normal"
}

foreach ($change in $changes) {
    $addr = "$($change.Col)$($change.Row)"
    Set-TextValue $ws.Range($addr) $change.Value
}

# Rows 46 and 47 swap places (USDe moves above Fetch.AI) with new volume data.
Set-TextValue $ws.Range("B46") "USDe"
Set-TextValue $ws.Range("C46") "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-TextValue $ws.Range("D46") "0.999"
Set-TextValue $ws.Range("E46") "  -0.04%  "

Set-TextValue $ws.Range("B47") "Fetch.AI"
Set-TextValue $ws.Range("C47") "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-TextValue $ws.Range("D47") "2.13"
Set-TextValue $ws.Range("E47") "  -6.39%  "
